# Append three new daily rows (2020-12-02 .. 2020-12-04) to the Indiana
# hospital ventilator report, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2020-12-02", 2184, 987, 750, 447, 2798, 399, 449, 1950, 45.19, 34.34, 20.47, 14.26, 16.05, 69.69),
    @("2020-12-03", 2185, 955, 783, 447, 2798, 383, 448, 1967, 43.71, 35.84, 20.46, 13.69, 16.01, 70.3),
    @("2020-12-04", 2182, 961, 773, 448, 2794, 368, 445, 1981, 44.04, 35.43, 20.53, 13.17, 15.93, 70.90000000000001)
)

$startRow = 281
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Column A holds a plain text date string (e.g. "2020-12-02"), not a
    # real date serial. Force text formatting before assigning so Excel's
    # autodetection doesn't convert it to a date value, then restore the
    # default "Normal" style so no stray per-cell formatting is left behind.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.Style = "Normal"

    for ($c = 1; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
